# Bump the header version from "(v2)" to "(v3)" and relocate Word's
# "_GoBack" (last-edit-location) bookmark to sit right after the edit,
# between the new "3" and the closing parenthesis - mirroring exactly
# what Word itself does when you retype a character in the document.

$d = $word.ActiveDocument

# Locate " (v2)" at the end of the title line so we work with absolute
# character offsets instead of hard-coding them.
$hit = $d.Content
$found = $hit.Find.Execute(" (v2)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find ' (v2)' in the document"
}
$openParenStart = $hit.Start      # start of " (v2)"
$digitStart = $openParenStart + 3 # position of the "2"  (' ','(','v','2',')')
$digitEnd = $digitStart + 1       # just after the "2"

# Pin the run boundary that precedes " (v2)" with a scratch bookmark so
# that the upcoming text edit cannot coalesce it into the neighbouring
# run - we delete this scratch bookmark again once the real edit is done.
$d.Bookmarks.Add("zzzScratchAnchor", $d.Range($openParenStart, $openParenStart)) | Out-Null

# "2" -> "3"
$digitRange = $d.Range($digitStart, $digitEnd)
$digitRange.Text = "3"

# Word drops "_GoBack" at the spot of the most recent edit - here, right
# after the newly typed "3" and before the closing ")". Re-adding a
# bookmark named "_GoBack" automatically relocates/replaces whichever
# "_GoBack" bookmark already existed elsewhere in the document.
$d.Bookmarks.Add("_GoBack", $d.Range($digitEnd, $digitEnd)) | Out-Null

# Drop the scratch bookmark now that the split is locked in.
$d.Bookmarks("zzzScratchAnchor").Delete()
